$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.974.55'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.901.01'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.94%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7432'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.57'
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.37'
$ws.Range("E9").Value = '  -6.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06885'
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08012'
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7524'
$ws.Range("E12").Value = '  -2.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.901.42'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.258'
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.36'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.137'
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.983.84'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.95'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007745'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.74'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.149.59'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.049'
$ws.Range("E24").Value = '  +6.70%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.294'
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.90'
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.77'
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1263'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.045'
$ws.Range("E29").Value = '  -4.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.349'
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.299'
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.034'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05378'
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.280'
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7368'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01942'
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.769'
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.205'
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4446'
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.42'
$ws.Range("E42").Value = '  -4.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.938'
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8319'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.705'
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.57'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.844'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.059.15'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.59'
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1162'
$ws.Range("E51").Value = '  -3.97%  '
Write-Host "Update complete"
